$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 2 data: A2 = 1, B2 = "Test"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Test"

# Update selection to B2
$ws.Range("B2").Select()
